$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.072.38"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "1.789.29"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.89"
$ws.Range("E5").Value = "  +1.33%  "
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.19"
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.295"
$ws.Range("E9").Value = "  +4.12%  "
$ws.Range("E10").Value = "  -2.20%  "
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("D12").Value = "2.047.44"
$ws.Range("E13").Value = "  +5.34%  "
$ws.Range("D14").Value = "1.789.85"
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.622"
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").Value = "34.073.02"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.04"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.73"
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.93"
$ws.Range("E21").Value = "  +2.81%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.09"
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("E24").Value = "  -2.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.10"
$ws.Range("E25").Value = "  +1.74%  "
$ws.Range("E26").Value = "  +2.58%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("E28").Value = "  +1.47%  "
$ws.Range("E30").Value = "  +1.76%  "
$ws.Range("E31").Value = "  +1.84%  "
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.60"
$ws.Range("E33").Value = "  +3.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.84"
$ws.Range("E34").Value = "  +1.62%  "
$ws.Range("D35").Value = "1.411.51"
$ws.Range("E35").Value = "  +1.53%  "
$ws.Range("E36").Value = "  +1.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0191"
$ws.Range("E37").Value = "  +3.10%  "
$ws.Range("E38").Value = "  +8.08%  "
$ws.Range("E39").Value = "  -0.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.67"
$ws.Range("E40").Value = "  +3.15%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.920"
$ws.Range("E42").Value = "  +1.32%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.35"
$ws.Range("E44").Value = "  +9.13%  "
$ws.Range("E45").Value = "  -4.76%  "
$ws.Range("E46").Value = "  +2.10%  "
$ws.Range("E47").Value = "  +3.40%  "
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "106.83"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").Value = "1.947.88"
$ws.Range("E50").Value = "  +0.45%  "
